# Add manual clusters analysis tabular results: extend the MSA clustering
# results sheet with a third clustering method ("Brute-force") contributing
# three new columns (H: ClusterSize, I: MinimumPercentIdentity, J: Average)
# alongside the existing CSANN (B:D) and CD-HIT (E:G) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (H1:J1), added as new shared strings ---
$ws.Range("H1").Value = "ClusterSize(Brute-force)"
$ws.Range("I1").Value = "MinimumPercentIdentity(Brute-force)"
$ws.Range("J1").Value = "Average(Brute-force)"

# Copy the header fill/format from the last existing header cell (G1) onto
# the three new header cells so they look consistent with B1:G1.
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# --- Brute-force clustering results for data rows 2-51 ---
# Each entry is (row, ClusterSize, MinimumPercentIdentity, Average)
$bruteForceData = @(
    @(2, 71, 83, 93.150499999999994),
    @(3, 67, 74, 90.077299999999994),
    @(4, 57, 86, 95.3352),
    @(5, 57, 77, 87.206100000000006),
    @(6, 53, 87, 95.493499999999997),
    @(7, 53, 82, 96.260499999999993),
    @(8, 51, 79, 88.7239),
    @(9, 45, 78, 92.509100000000004),
    @(10, 38, 95, 97.597399999999993),
    @(11, 37, 69, 93.572100000000006),
    @(12, 36, 83, 94.225399999999993),
    @(13, 36, 80, 89.001599999999996),
    @(14, 36, 81, 90.963499999999996),
    @(15, 35, 78, 90.253799999999998),
    @(16, 34, 79, 92.183599999999998),
    @(17, 34, 83, 95.739800000000002),
    @(18, 34, 77, 85.138999999999996),
    @(19, 31, 76, 87.563400000000001),
    @(20, 30, 72, 87.089699999999993),
    @(21, 28, 85, 90.454999999999998),
    @(22, 27, 77, 90.6952),
    @(23, 25, 75, 85.423299999999998),
    @(24, 25, 78, 88.713300000000004),
    @(25, 24, 78, 88.724599999999995),
    @(26, 24, 70, 91.130399999999995),
    @(27, 24, 75, 82.837000000000003),
    @(28, 24, 79, 90.884100000000004),
    @(29, 23, 73, 82.533600000000007),
    @(30, 23, 76, 83.905100000000004),
    @(31, 23, 75, 85.454499999999996),
    @(32, 22, 79, 88.584400000000002),
    @(33, 22, 76, 84.116900000000001),
    @(34, 22, 75, 86.515199999999993),
    @(35, 22, 94, 96.29),
    @(36, 22, 79, 88.857100000000003),
    @(37, 22, 79, 91.142899999999997),
    @(38, 21, 75, 88.5),
    @(39, 20, 81, 90.710499999999996),
    @(40, 20, 78, 87.042100000000005),
    @(41, 20, 71, 83.284199999999998),
    @(42, 20, 79, 89.868399999999994),
    @(43, 20, 70, 80.126300000000001),
    @(44, 20, 75, 83.1053),
    @(45, 20, 91, 95.547399999999996),
    @(46, 20, 83, 90.257900000000006),
    @(47, 19, 78, 85.795299999999997),
    @(48, 19, 85, 90.532200000000003),
    @(49, 19, 82, 89.707599999999999),
    @(50, 19, 78, 91.157899999999998),
    @(51, 19, 89, 94.251499999999993)
)

foreach ($entry in $bruteForceData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 8).Value  = $entry[1]
    $ws.Cells.Item($r, 9).Value  = $entry[2]
    $ws.Cells.Item($r, 10).Value = $entry[3]
}

# --- Column widths: widen A:G slightly and size the new H:I:J columns
#     (closest values achievable through the ColumnWidth property, which
#     Excel stores/quantizes internally before writing the sheet XML). ---
$ws.Columns.Item(1).ColumnWidth  = 2.6666666666666665
$ws.Columns.Item(2).ColumnWidth  = 18
$ws.Columns.Item(3).ColumnWidth  = 31
$ws.Columns.Item(4).ColumnWidth  = 15.333333333333334
$ws.Columns.Item(5).ColumnWidth  = 17.666666666666668
$ws.Columns.Item(6).ColumnWidth  = 30.666666666666668
$ws.Columns.Item(7).ColumnWidth  = 15
$ws.Columns.Item(8).ColumnWidth  = 22
$ws.Columns.Item(9).ColumnWidth  = 35
$ws.Columns.Item(10).ColumnWidth = 19.333333333333332

# --- Selection now covers the full extended table ---
$ws.Range("A2:J51").Select()

# --- Page setup: explicit portrait orientation ---
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait
